# Add the new "Unique Data" sheet after the existing "ManageNetwork" sheet
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "Unique Data"

# Populate the new sheet
$ws2.Range("A1").Value = "Columns for which unique data to be provided"
$ws2.Range("B1").Value = "Comments"
$ws2.Range("A2").Value = "DeviceName"
$ws2.Range("B2").Value = "Based on the device name, all other column data should be provided from Application"

# Update selection on ManageNetwork sheet
$ws1.Range("B3").Select()
